$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the price column header label from "价格" to "进院价格"
$ws.Range("K1").Value = "进院价格"

# Match the resulting selection change observed after the edit
$ws.Range("A3").Select()
